# Update "想去人数" (F column) values across the four sheets to reflect
# the latest scrape output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 845
$ws1.Range("F3").Value = 13712
$ws1.Range("F4").Value = 13493
$ws1.Range("F5").Value = 1047
$ws1.Range("F9").Value = 79
$ws1.Range("F10").Value = 22
$ws1.Range("F12").Value = 751
$ws1.Range("F13").Value = 2136
$ws1.Range("F16").Value = 70
$ws1.Range("F17").Value = 112
$ws1.Range("F19").Value = 511
$ws1.Range("F21").Value = 382
$ws1.Range("F22").Value = 314
$ws1.Range("F23").Value = 155
$ws1.Range("F24").Value = 821
$ws1.Range("F25").Value = 73

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 70
$ws2.Range("F7").Value = 1443
$ws2.Range("F10").Value = 59
$ws2.Range("F12").Value = 8

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 218

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 218
$ws4.Range("F3").Value = 845
$ws4.Range("F4").Value = 13712
$ws4.Range("F5").Value = 13493
$ws4.Range("F6").Value = 1047
$ws4.Range("F10").Value = 79
$ws4.Range("F11").Value = 22
$ws4.Range("F13").Value = 751
$ws4.Range("F16").Value = 2136
$ws4.Range("F19").Value = 70
$ws4.Range("F20").Value = 112
$ws4.Range("F23").Value = 70
$ws4.Range("F26").Value = 511
$ws4.Range("F28").Value = 382
$ws4.Range("F29").Value = 314
$ws4.Range("F30").Value = 156
$ws4.Range("F31").Value = 821
$ws4.Range("F33").Value = 1443
$ws4.Range("F36").Value = 73
$ws4.Range("F37").Value = 59
$ws4.Range("F40").Value = 8
